$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "255.54"
Set-TextValue "E2" "4.16%"
Set-TextValue "G2" "5"

Set-TextValue "D3" "28.00"
Set-TextValue "E3" "-4.78%"
Set-TextValue "G3" "5"

Set-TextValue "D4" "5.374"
Set-TextValue "E4" "4.37%"
Set-TextValue "G4" "5"

Set-TextValue "D5" "0.05826"
Set-TextValue "E5" "0.85%"
Set-TextValue "G5" "5"

Set-TextValue "D6" "6.719"
Set-TextValue "E6" "1.36%"
Set-TextValue "G6" "5"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D7" "0.8699"
Set-TextValue "E7" "1.57%"
Set-TextValue "G7" "5"

$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D8" "0.9122"
Set-TextValue "E8" "6.51%"
Set-TextValue "G8" "5"

$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1417"
Set-TextValue "E9" "3.99%"
Set-TextValue "G9" "5"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07229"
Set-TextValue "E10" "2.84%"
Set-TextValue "G10" "5"

$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.03179"
Set-TextValue "E11" "1.05%"
Set-TextValue "G11" "5"

$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D12" "0.09242"
Set-TextValue "E12" "-1.34%"
Set-TextValue "G12" "5"

$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D13" "0.001544"
Set-TextValue "E13" "1.47%"
Set-TextValue "G13" "5"

$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D14" "0.01058"
Set-TextValue "E14" "1,663.82%"
Set-TextValue "G14" "5"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.006045"
Set-TextValue "E15" "0.54%"
Set-TextValue "G15" "5"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.493"
Set-TextValue "E16" "0.22%"
Set-TextValue "G16" "5"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "3.225"
Set-TextValue "E17" "1.54%"
Set-TextValue "G17" "5"

Set-TextValue "E18" "5.07%"
Set-TextValue "G18" "5"

Set-TextValue "D19" "0.3168"
Set-TextValue "E19" "-1.00%"
Set-TextValue "G19" "5"

Set-TextValue "D20" "0.03459"
Set-TextValue "E20" "4.82%"
Set-TextValue "G20" "5"

Set-TextValue "G21" "5"

Set-TextValue "D22" "3.542"
Set-TextValue "E22" "6.78%"
Set-TextValue "G22" "5"

Set-TextValue "D23" "0.04171"
Set-TextValue "E23" "0.93%"
Set-TextValue "G23" "5"

Set-TextValue "D24" "0.1379"
Set-TextValue "E24" "-1.49%"
Set-TextValue "G24" "5"

Set-TextValue "D25" "0.001227"
Set-TextValue "E25" "0.05%"
Set-TextValue "G25" "5"

Set-TextValue "D26" "0.004867"
Set-TextValue "E26" "17.93%"
Set-TextValue "G26" "5"

Set-TextValue "D27" "0.0001200"
Set-TextValue "E27" "-0.82%"
Set-TextValue "G27" "5"

Set-TextValue "D28" "0.0001455"
Set-TextValue "E28" "0.70%"
Set-TextValue "G28" "5"

Set-TextValue "G29" "5"

Set-TextValue "G30" "5"

Set-TextValue "G31" "5"

Set-TextValue "G32" "5"

Set-TextValue "G33" "5"

Set-TextValue "G34" "5"

Set-TextValue "G35" "5"

Set-TextValue "G36" "5"

Set-TextValue "G37" "5"

Set-TextValue "G38" "5"

Set-TextValue "G39" "5"

Set-TextValue "D40" "0.03851"
Set-TextValue "E40" "3.35%"
Set-TextValue "G40" "5"

Set-TextValue "D41" "0.005754"
Set-TextValue "E41" "-1.96%"
Set-TextValue "G41" "5"

Set-TextValue "D42" "0.1100"
Set-TextValue "E42" "3.03%"
Set-TextValue "G42" "5"

Set-TextValue "D43" "0.002200"
Set-TextValue "E43" "0.01%"
Set-TextValue "G43" "5"

Set-TextValue "D44" "0.009948"
Set-TextValue "E44" "8.44%"
Set-TextValue "G44" "5"

Set-TextValue "D45" "0.00005277"
Set-TextValue "E45" "-0.06%"
Set-TextValue "G45" "5"

Set-TextValue "D46" "0.00000000750"
Set-TextValue "E46" "0.01%"
Set-TextValue "G46" "5"

Set-TextValue "D47" "0.1000"
Set-TextValue "E47" "72.43%"
Set-TextValue "G47" "5"

Set-TextValue "D48" "0.002199"
Set-TextValue "E48" "1.19%"
Set-TextValue "G48" "5"

Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "0.01%"
Set-TextValue "G49" "5"

Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "0.01%"
Set-TextValue "G50" "5"

Set-TextValue "G51" "5"
